$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.23"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.639.53"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0610"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.873.15"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.649.90"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "27.650.16"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "1.447.39"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.916"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.70%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "1.782.68"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "
